$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend formatting down by copying the prior two data rows (54:55) into the
# new rows (56:57) so the new rows inherit the same cell styles (date /
# decimal / integer number formats) instead of creating new style entries.
$ws.Range("A54:F55").Copy()
$ws.Range("A56:F57").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row 56: 四方坪站
$ws.Range("A56").Value = 45928
$ws.Range("B56").Value = "四方坪站"
$ws.Range("C56").Value = 9093.57
$ws.Range("D56").Value = 7618.78
$ws.Range("E56").Value = 3124.72
$ws.Range("F56").Value = 386

# New row 57: 高岭站
$ws.Range("A57").Value = 45928
$ws.Range("B57").Value = "高岭站"
$ws.Range("C57").Value = 5876.99
$ws.Range("D57").Value = 4743.26
$ws.Range("E57").Value = 1432.67
$ws.Range("F57").Value = 206

$ws.Range("I56").Select()
